$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordered list of requirement IDs that should occupy rows 2..27 (column A).
$values = @(
    "REQ0782445",
    "REQ0782366",
    "REQ0782025",
    "REQ0781660",
    "REQ0781412",
    "REQ0781385",
    "REQ0780724",
    "REQ0780718",
    "REQ0780623",
    "REQ0777639",
    "REQ0770139",
    "REQ0767628",
    "REQ0767543",
    "REQ0767475",
    "REQ0767407",
    "REQ0767333",
    "REQ0767300",
    "REQ0767209",
    "REQ0767158",
    "REQ0759241",
    "REQ0759224",
    "REQ0759215",
    "REQ0759187",
    "REQ0759135",
    "REQ0758818",
    "REQ0753679"
)

$startRow = 2
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
